$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 331
$ws.Range("F4").Value = 644
$ws.Range("F7").Value = 2149
$ws.Range("F8").Value = 892
$ws.Range("F9").Value = 847
$ws.Range("F10").Value = 411
$ws.Range("F11").Value = 87
$ws.Range("F12").Value = 430
$ws.Range("F13").Value = 322
$ws.Range("F14").Value = 99
$ws.Range("F15").Value = 895
$ws.Range("F17").Value = 35
$ws.Range("F18").Value = 1762
$ws.Range("F22").Value = 59
$ws.Range("F24").Value = 1456
$ws.Range("F26").Value = 528
$ws.Range("F27").Value = 352
$ws.Range("F28").Value = 610
$ws.Range("F29").Value = 420
$ws.Range("F30").Value = 2428
$ws.Range("F31").Value = 382
$ws.Range("F32").Value = 93
$ws.Range("G32").Value = 88
$ws.Range("F33").Value = 169
$ws.Range("F34").Value = 601
$ws.Range("F35").Value = 477
$ws.Range("F36").Value = 193
$ws.Range("F37").Value = 922
$ws.Range("F38").Value = 704
$ws.Range("F40").Value = 487
$ws.Range("F41").Value = 466

# Sheet: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F11").Value = 55
$ws.Range("F15").Value = 4
$ws.Range("F22").Value = 119
$ws.Range("F23").Value = 108
$ws.Range("F24").Value = 435

# Sheet: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 228
$ws.Range("F3").Value = 2921
$ws.Range("F5").Value = 245
$ws.Range("F6").Value = 318

# Sheet: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 228
$ws.Range("F6").Value = 245
$ws.Range("F7").Value = 644
$ws.Range("F10").Value = 892
$ws.Range("F11").Value = 847
$ws.Range("F12").Value = 411
$ws.Range("F13").Value = 87
$ws.Range("F14").Value = 322
$ws.Range("F15").Value = 99
$ws.Range("F17").Value = 895
$ws.Range("F20").Value = 35
$ws.Range("F21").Value = 318
$ws.Range("F22").Value = 1762
$ws.Range("F26").Value = 55
$ws.Range("F30").Value = 1456
$ws.Range("F31").Value = 4
$ws.Range("F33").Value = 528
$ws.Range("F34").Value = 352
$ws.Range("F35").Value = 610
$ws.Range("F36").Value = 420
$ws.Range("F37").Value = 93
$ws.Range("G37").Value = 88
$ws.Range("F38").Value = 169
$ws.Range("F39").Value = 477
$ws.Range("F40").Value = 193
$ws.Range("F41").Value = 922
$ws.Range("F44").Value = 108
$ws.Range("F45").Value = 435
$ws.Range("F46").Value = 704
$ws.Range("F48").Value = 487
$ws.Range("F49").Value = 466
